$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Fri. 11/21" heading paragraph: drop the stray
#    <w:rFonts w:hint="eastAsia"/> that lived in the paragraph
#    mark's own rPr (w:pPr/w:rPr). The run-level rPr keeps its hint.
# ------------------------------------------------------------------
$p = $d.Paragraphs(2)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:u w:val="single"/></w:rPr><w:t>Fri. 11/21</w:t></w:r></w:p>
'@)

# ------------------------------------------------------------------
# 2) "Sat. 11/22" heading paragraph: same cleanup.
# ------------------------------------------------------------------
$p = $d.Paragraphs(10)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:u w:val="single"/></w:rPr><w:t>Sat. 11/22</w:t></w:r></w:p>
'@)

# ------------------------------------------------------------------
# 3) "Comment each Code" paragraph becomes a 2-line-long cache
#    testing note split across four runs (with a spell-checked
#    "MainMemory" run), and loses its w:pPr entirely.
# ------------------------------------------------------------------
$p = $d.Paragraphs(11)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>DONE-</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Comment </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">L1Cache and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>MainMemory</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@)

# ------------------------------------------------------------------
# 4) "Read cache line from DRAM" list item: drop the redundant
#    <w:rFonts w:hint="eastAsia"/> sitting in the paragraph mark's
#    rPr inside w:pPr (the numPr/pStyle stay untouched).
# ------------------------------------------------------------------
$p = $d.Paragraphs(19)
$p.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Read cache line from DRAM</w:t></w:r></w:p>
'@)

Write-Output "edits applied"
